# Add files via upload
# - Re-bucket the A/C columns (range + name) for rows 2-5, adding two new
#   "Neutral"-styled rows (200-399/400-449) and two new plain rows
#   (100-199/450-499) that sit alongside the existing per-ID rows in column B.
# - Column B (the 01..17 id list) is untouched in value/style except that the
#   two new A/C pairs line up against rows 4 and 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: now "200-399" / 02 / "Andrew", styled like the "Good" banding
#     used elsewhere in the sheet, but with the new "Neutral" (yellow) look.
$ws.Range("A2").Value = "200-399"
$ws.Range("B2").Value = "02"
$ws.Range("C2").Value = "Andrew"
$ws.Range("A2:C2").Style = "Neutral"
$ws.Range("B2").NumberFormat = "@"

# --- Row 3: "400-449" / 03 / "Miho"
$ws.Range("A3").Value = "400-449"
$ws.Range("B3").Value = "03"
$ws.Range("C3").Value = "Miho"
$ws.Range("A3:C3").Style = "Neutral"
$ws.Range("B3").NumberFormat = "@"

# --- Row 4: new A4/C4 values added next to the existing B4 ("04"); no
#     special fill, matching the rest of the plain (unstyled) column A/C.
$ws.Range("A4").Value = "100-199"
$ws.Range("C4").Value = "Miho"

# --- Row 5: new A5/C5 values added next to the existing B5 ("05")
$ws.Range("A5").Value = "450-499"
$ws.Range("C5").Value = "Andrew"

# --- Selection moved as part of the author's last interaction with the sheet
[void]$ws.Range("E14").Select()
